{"js": "const replacements = [\n  [\"2024-09-09 Monday\", \"2024-09-10 Tuesday\"],\n  [\"432\\u00d79=3888\", \"668\\u00d76=4008\"],\n  [\"108\\u00d78=864\", \"189\\u00d75=945\"],\n  [\"937\\u00d73=2811\", \"259\\u00d79=2331\"],\n  [\"597\\u00d73=1791\", \"208\\u00d77=1456\"],\n  [\"863\\u00d78=6904\", \"529\\u00d79=4761\"],\n  [\"829\\u00d74=3316\", \"914\\u00d74=3656\"],\n  [\"195\\u00d72=390\", \"858\\u00d74=3432\"],\n  [\"258\\u00d75=1290\", \"485\\u00d76=2910\"],\n  [\"504\\u00d77=3528\", \"161\\u00d74=644\"],\n  [\"103\\u00d79=927\", \"541\\u00d75=2705\"],\n  [\"908\\u00d74=3632\", \"887\\u00d79=7983\"],\n  [\"877\\u00d72=1754\", \"729\\u00d76=4374\"],\n  [\"170\\u00d72=340\", \"695\\u00d76=4170\"],\n  [\"361\\u00d79=3249\", \"943\\u00d77=6601\"],\n  [\"815\\u00d74=3260\", \"874\\u00d76=5244\"],\n  [\"197\\u00d79=1773\", \"179\\u00d78=1432\"],\n  [\"938\\u00d75=4690\", \"274\\u00d73=822\"],\n  [\"703\\u00d75=3515\", \"883\\u00d76=5298\"],\n  [\"514\\u00d73=1542\", \"829\\u00d73=2487\"],\n  [\"157\\u00d76=942\", \"965\\u00d77=6755\"],\n  [\"509\\u00d72=1018\", \"992\\u00d75=4960\"],\n  [\"622\\u00d72=1244\", \"182\\u00d72=364\"],\n  [\"829\\u00d72=1658\", \"712\\u00d75=3560\"],\n  [\"383\\u00d77=2681\", \"954\\u00d74=3816\"],\n  [\"837\\u00d77=5859\", \"566\\u00d76=3396\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-09 Monday\", \"2024-09-10 Tuesday\"),\n    @(\"432\u00d79=3888\", \"668\u00d76=4008\"),\n    @(\"108\u00d78=864\", \"189\u00d75=945\"),\n    @(\"937\u00d73=2811\", \"259\u00d79=2331\"),\n    @(\"597\u00d73=1791\", \"208\u00d77=1456\"),\n    @(\"863\u00d78=6904\", \"529\u00d79=4761\"),\n    @(\"829\u00d74=3316\", \"914\u00d74=3656\"),\n    @(\"195\u00d72=390\", \"858\u00d74=3432\"),\n    @(\"258\u00d75=1290\", \"485\u00d76=2910\"),\n    @(\"504\u00d77=3528\", \"161\u00d74=644\"),\n    @(\"103\u00d79=927\", \"541\u00d75=2705\"),\n    @(\"908\u00d74=3632\", \"887\u00d79=7983\"),\n    @(\"877\u00d72=1754\", \"729\u00d76=4374\"),\n    @(\"170\u00d72=340\", \"695\u00d76=4170\"),\n    @(\"361\u00d79=3249\", \"943\u00d77=6601\"),\n    @(\"815\u00d74=3260\", \"874\u00d76=5244\"),\n    @(\"197\u00d79=1773\", \"179\u00d78=1432\"),\n    @(\"938\u00d75=4690\", \"274\u00d73=822\"),\n    @(\"703\u00d75=3515\", \"883\u00d76=5298\"),\n    @(\"514\u00d73=1542\", \"829\u00d73=2487\"),\n    @(\"157\u00d76=942\", \"965\u00d77=6755\"),\n    @(\"509\u00d72=1018\", \"992\u00d75=4960\"),\n    @(\"622\u00d72=1244\", \"182\u00d72=364\"),\n    @(\"829\u00d72=1658\", \"712\u00d75=3560\"),\n    @(\"383\u00d77=2681\", \"954\u00d74=3816\"),\n    @(\"837\u00d77=5859\", \"566\u00d76=3396\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
